$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 665.55554
$ws.Range("I4").Value = 584.2857
$ws.Range("J4").Value = 950
$ws.Range("K4").Value = 584.2857
$ws.Range("L4").Value = 950
$ws.Range("M4").Value = -470.2857
$ws.Range("N4").Value = -1178
$ws.Range("H9").Value = 85
$ws.Range("I9").Value = 85
$ws.Range("K9").Value = 85
$ws.Range("M9").Value = 84
$ws.Range("H40").Value = 2123.0625
$ws.Range("I40").Value = 1924.1428
$ws.Range("J40").Value = 2277.7778
$ws.Range("K40").Value = 1924.1428
$ws.Range("L40").Value = 2277.7778
$ws.Range("M40").Value = -1749.1428
$ws.Range("N40").Value = -2627.7778
$ws.Range("H131").Value = 2438.8
$ws.Range("I131").Value = 1931.3334
$ws.Range("J131").Value = 3200
$ws.Range("K131").Value = 5794.0002
$ws.Range("L131").Value = 9600
$ws.Range("M131").Value = -754.0002000000004
$ws.Range("N131").Value = -19680
$ws.Range("H137").Value = 2872.8823
$ws.Range("I137").Value = 2144.8333
$ws.Range("J137").Value = 4620.2
$ws.Range("K137").Value = 6434.499899999999
$ws.Range("L137").Value = 13860.6
$ws.Range("M137").Value = -3884.499899999999
$ws.Range("N137").Value = -18960.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5475.2
$ws.Range("I32").Value = 3754.5454
$ws.Range("K32").Value = 3754.5454
$ws.Range("M32").Value = -3467.5454
$ws.Range("H45").Value = 1770.1818
$ws.Range("I45").Value = 1770.8
$ws.Range("K45").Value = 1770.8
$ws.Range("M45").Value = -1393.8
$ws.Range("H61").Value = 1948.8695
$ws.Range("I61").Value = 1748.0526
$ws.Range("J61").Value = 2902.75
$ws.Range("K61").Value = 1748.0526
$ws.Range("L61").Value = 2902.75
$ws.Range("M61").Value = -1536.0526
$ws.Range("N61").Value = -3326.75
$ws.Range("H102").Value = 1225
$ws.Range("I102").Value = 1225
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1225
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 397
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 5956.5713
$ws.Range("I110").Value = 6449.3335
$ws.Range("K110").Value = 6449.3335
$ws.Range("M110").Value = -4404.3335
$ws.Range("H132").Value = 2047
$ws.Range("J132").Value = 2048.8572
$ws.Range("L132").Value = 6146.571599999999
$ws.Range("N132").Value = -11206.5716
$ws.Range("H136").Value = 1948.8695
$ws.Range("I136").Value = 1748.0526
$ws.Range("J136").Value = 2902.75
$ws.Range("K136").Value = 5244.1578
$ws.Range("L136").Value = 8708.25
$ws.Range("M136").Value = -2694.1578
$ws.Range("N136").Value = -13808.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 25000
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H134").Value = 2395.6086
$ws.Range("I134").Value = 2113.25
$ws.Range("K134").Value = 6339.75
$ws.Range("M134").Value = -3804.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10002.667
$ws.Range("J3").Value = 10002.667
$ws.Range("L3").Value = 10002.667
$ws.Range("N3").Value = -10228.667
$ws.Range("H22").Value = 313.55554
$ws.Range("I22").Value = 277.75
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 277.75
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 72.25
$ws.Range("N22").Value = -1300
$ws.Range("H58").Value = 2729.054
$ws.Range("I58").Value = 1253.619
$ws.Range("K58").Value = 1253.619
$ws.Range("M58").Value = -1050.619
$ws.Range("H60").Value = 12217.5
$ws.Range("I60").Value = 9698.733
$ws.Range("J60").Value = 49999
$ws.Range("K60").Value = 9698.733
$ws.Range("L60").Value = 49999
$ws.Range("M60").Value = -9187.733
$ws.Range("N60").Value = -51021
$ws.Range("H122").Value = 2311.0588
$ws.Range("I122").Value = 2268.0625
$ws.Range("K122").Value = 6804.1875
$ws.Range("M122").Value = -4354.1875
$ws.Range("H132").Value = 1751.5344
$ws.Range("I132").Value = 1465.3273
$ws.Range("J132").Value = 6998.6665
$ws.Range("K132").Value = 4395.9819
$ws.Range("L132").Value = 20995.9995
$ws.Range("M132").Value = -1865.9819
$ws.Range("N132").Value = -26055.9995
$ws.Range("H133").Value = 34994.5
$ws.Range("I133").Value = 34994.5
$ws.Range("K133").Value = 34994.5
$ws.Range("M133").Value = -32464.5
$ws.Range("H136").Value = 2729.054
$ws.Range("I136").Value = 1253.619
$ws.Range("K136").Value = 3760.857
$ws.Range("M136").Value = -1210.857
$ws.Range("H141").Value = 128002.29
$ws.Range("J141").Value = 149684
$ws.Range("L141").Value = 149684
$ws.Range("N141").Value = -160044

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 412.5926
$ws.Range("I107").Value = 276.14285
$ws.Range("J107").Value = 460.35
$ws.Range("K107").Value = 828.4285500000001
$ws.Range("L107").Value = 1381.05
$ws.Range("M107").Value = 1091.57145
$ws.Range("N107").Value = -5221.05
$ws.Range("H140").Value = 1418.0834
$ws.Range("I140").Value = 1244.909
$ws.Range("K140").Value = 3734.727
$ws.Range("M140").Value = 1445.273

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 267160.5
$ws.Range("J11").Value = 345130.66
$ws.Range("L11").Value = 345130.66
$ws.Range("N11").Value = -345408.66
$ws.Range("H29").Value = 726.25
$ws.Range("I29").Value = 162
$ws.Range("J29").Value = 1666.6666
$ws.Range("K29").Value = 162
$ws.Range("L29").Value = 1666.6666
$ws.Range("M29").Value = 128
$ws.Range("N29").Value = -2246.6666
$ws.Range("H107").Value = 112.333336
$ws.Range("I107").Value = 111.8
$ws.Range("J107").Value = 115
$ws.Range("K107").Value = 111.8
$ws.Range("L107").Value = 115
$ws.Range("M107").Value = 1808.2
$ws.Range("N107").Value = -3955

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4666
$ws.Range("I7").Value = 4666
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4666
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4554
$ws.Range("N7").ClearContents()
$ws.Range("H46").Value = 2772.36
$ws.Range("I46").Value = 2045.0555
$ws.Range("J46").Value = 4642.5713
$ws.Range("K46").Value = 2045.0555
$ws.Range("L46").Value = 4642.5713
$ws.Range("M46").Value = -1857.0555
$ws.Range("N46").Value = -5018.5713
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1798
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 4666
$ws.Range("I126").Value = 4666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11528
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 4525.857
$ws.Range("I136").Value = 3973
$ws.Range("K136").Value = 11919
$ws.Range("M136").Value = -9369

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 352964
$ws.Range("I3").Value = 514994
$ws.Range("J3").Value = 28904
$ws.Range("K3").Value = 514994
$ws.Range("L3").Value = 28904
$ws.Range("M3").Value = -514880
$ws.Range("N3").Value = -29132
$ws.Range("H21").Value = 18486.166
$ws.Range("J21").Value = 18486.166
$ws.Range("L21").Value = 18486.166
$ws.Range("N21").Value = -18956.166
$ws.Range("H35").Value = 18486.166
$ws.Range("J35").Value = 18486.166
$ws.Range("L35").Value = 18486.166
$ws.Range("N35").Value = -19066.166
$ws.Range("H64").Value = 66500
$ws.Range("J64").Value = 66500
$ws.Range("L64").Value = 66500
$ws.Range("N64").Value = -66996
$ws.Range("H67").Value = 66500
$ws.Range("J67").Value = 66500
$ws.Range("L67").Value = 66500
$ws.Range("N67").Value = -68216
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H122").Value = 3271.3333
$ws.Range("I122").Value = 3555.25
$ws.Range("K122").Value = 10665.75
$ws.Range("M122").Value = -8215.75
$ws.Range("H124").Value = 19966
$ws.Range("J124").Value = 19966
$ws.Range("L124").Value = 19966
$ws.Range("N124").Value = -29786
$ws.Range("H126").Value = 2288.3333
$ws.Range("I126").Value = 1737.5
$ws.Range("K126").Value = 5212.5
$ws.Range("M126").Value = -2742.5
$ws.Range("H132").Value = 32510.934
$ws.Range("I132").Value = 40247.668
$ws.Range("K132").Value = 120743.004
$ws.Range("M132").Value = -118213.004
$ws.Range("H136").Value = 1629.8334
$ws.Range("I136").Value = 1575.4117
$ws.Range("J136").Value = 2555
$ws.Range("K136").Value = 4726.2351
$ws.Range("L136").Value = 7665
$ws.Range("M136").Value = -2176.2351
$ws.Range("N136").Value = -12765

